$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '45.458.44'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -6.84%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.681.34'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +0.17%  '

$ws.Range("E4").Value = '  +0.31%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '309.73'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.75%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '98.60'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -7.58%  '

$ws.Range("E7").Value = '  -3.15%  '

$ws.Range("E8").Value = '  +0.15%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.582'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -3.53%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '38.47'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  -7.11%  '

$ws.Range("E11").Value = '  -3.04%  '

$ws.Range("E12").Value = '  -4.58%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '3.089.96'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.27%  '

$ws.Range("E14").Value = '  +0.17%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.690.09'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  +0.71%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.930'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.49%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '15.16'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.84%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '45.428.26'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -6.74%  '

$ws.Range("E19").Value = '  -2.81%  '

$ws.Range("E20").Value = '  -0.76%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.84'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  -4.94%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '75.02'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +1.35%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '281.45'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -0.77%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '3.04'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.87%  '

$ws.Range("B25").Value = 'EthereumClassic'
$ws.Range("C25").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '31.21'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +1.99%  '

$ws.Range("B26").Value = 'ImmutableX'
$ws.Range("C26").Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.24'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.96%  '

$ws.Range("E27").Value = '  -0.20%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '4.04'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -2.20%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '10.64'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.11%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.23'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -4.10%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '38.34'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -7.22%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '6.22'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.71%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.79'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.56%  '

$ws.Range("E34").Value = '  +2.85%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '156.08'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +1.49%  '

$ws.Range("E36").Value = '  -3.13%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '2.83'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.43%  '

$ws.Range("E38").Value = '  -4.84%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '25.73'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +11.41%  '

$ws.Range("E40").Value = '  -1.22%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '16.38'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -2.52%  '

$ws.Range("E42").Value = '  -4.39%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0326'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -4.66%  '

$ws.Range("E44").Value = '  -8.79%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.116.05'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -4.75%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.998'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +0.30%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '93.56'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -6.00%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '111.91'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -3.68%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '9.30'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -7.74%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '2.941.39'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.38%  '

$ws.Range("E51").Value = '  -3.32%  '
